$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (ECs / Vip / Sctr / Inflammatory-Mac) ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7349876666666667
$ws.Range("H2").Value = 2.204963
$ws.Range("M2").Value = 0.02521133333333334
$ws.Range("N2").Value = 0.07563400000000001
$ws.Range("O2").Value = 0.2787497281937693
$ws.Range("P2").Value = 0.2787497281937693
$ws.Range("Q2").Value = 0.01853001906022222
$ws.Range("R2").Value = 0.166770171542
$ws.Range("S2").Value = 0.2787497281937693
$ws.Range("T2").Value = 0.2787497281937693

# --- Update row 3: target cluster changes from "MuSCs" to "Resolving-Mac", with new values ---
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7349876666666667
$ws.Range("H3").Value = 2.204963
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.065233
$ws.Range("N3").Value = 0.195699
$ws.Range("O3").Value = 0.7212502718062307
$ws.Range("P3").Value = 0.7212502718062307
$ws.Range("Q3").Value = 0.04794545045966667
$ws.Range("R3").Value = 0.4315090541370001
$ws.Range("S3").Value = 0.7212502718062307
$ws.Range("T3").Value = 0.7212502718062307

# --- Remove old row 4 (was "Resolving-Mac", now merged into row 3) ---
$ws.Rows.Item(4).Delete()
